# edit.ps1 - applies the draft-gandhi-spring-stamp-srpm-05.pptx revision
#
# Helper: replace a paragraph's visible text without letting the host's
# text-diff/run-splitting logic fragment a single run into several runs
# (Paragraphs(i,1).Text = "..." tends to split on common substrings, while
# re-targeting the same span through .Characters(start,len) keeps one run).
function Set-ParaText {
    param($TextRange, [int]$Index, [string]$NewText)

    $para = $TextRange.Paragraphs($Index, 1)
    $start = $para.Start
    $len = $para.Length
    $span = $TextRange.Characters($start, $len)
    $span.Text = $NewText
}

# Helper: set the font size of a whole paragraph's run(s) by paragraph index.
function Set-ParaFontSize {
    param($TextRange, [int]$Index, $Size)

    $para = $TextRange.Paragraphs($Index, 1)
    $para.Font.Size = $Size
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 3 - "Requirements and Scope"
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(2).TextFrame.TextRange
Set-ParaText $tr3 2 "Delay and Loss Measurement "
Set-ParaText $tr3 3 "Links and end-to-end P2P/P2MP SR paths"

# ---------------------------------------------------------------------------
# Slide 5 - "Session-Sender Test Packet for Links"
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$tr5 = $s5.Shapes.Item(3).TextFrame.TextRange
Set-ParaText $tr5 3 "TTL is set 1."

# ---------------------------------------------------------------------------
# Slide 6 - "Session-Sender Test Packet for SR-MPLS and SRv6 Policy"
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$tr6 = $s6.Shapes.Item(4).TextFrame.TextRange
Set-ParaText $tr6 8 "Color only SR-MPLS Policy:"
Set-ParaText $tr6 11 "P2MP SR-MPLS Policy"

# ---------------------------------------------------------------------------
# Slide 7 - "  Session-Reflector Test Packet"
# ---------------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$tr7a = $s7.Shapes.Item(2).TextFrame.TextRange
Set-ParaText $tr7a 2 "Link - Use Control Code from the received test packet if set."
Set-ParaText $tr7a 3 "SR path - Use Segment List from Return Path TLV if present in received test packet."

$tr7b = $s7.Shapes.Item(4).TextFrame.TextRange
Set-ParaText $tr7b 17 "               Figure 4: Session-Reflector Test Packet"

# ---------------------------------------------------------------------------
# Slide 8 - "ECMP Support for SR Path"
#   - move/resize the content placeholder slightly
#   - shrink every bullet from 18pt to 16pt (text itself is unchanged)
# ---------------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$sh8 = $s8.Shapes.Item(2)
$sh8.Left = 36
$sh8.Top = 74.16339
$tr8 = $sh8.TextFrame.TextRange
for ($i = 1; $i -le 8; $i++) {
    Set-ParaFontSize $tr8 $i 16
}

# ---------------------------------------------------------------------------
# Slide 9 - "Performance Measurement Modes"
#   - every bullet shrinks from 20pt to 18pt
#   - several bullets get re-worded
# ---------------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$tr9 = $s9.Shapes.Item(3).TextFrame.TextRange

Set-ParaText $tr9 1 "One-way Delay Measurement Mode"
Set-ParaFontSize $tr9 1 18

Set-ParaText $tr9 2 "Existing default behavior"
Set-ParaFontSize $tr9 2 18

Set-ParaText $tr9 3 "Two-way Delay Measurement Mode"
Set-ParaFontSize $tr9 3 18

Set-ParaFontSize $tr9 4 18

Set-ParaText $tr9 5 "Link - Use Control Code from the received test packet"
Set-ParaFontSize $tr9 5 18

Set-ParaText $tr9 6 "SR path - Use Return Path TLV for STAMP from the received test packet"
Set-ParaFontSize $tr9 6 18

Set-ParaFontSize $tr9 7 18

Set-ParaText $tr9 8 "  Session-sender test packet carries the return path in the header"
Set-ParaFontSize $tr9 8 18

# ---------------------------------------------------------------------------
# Slide 10 - "Example PM Metrics"
# ---------------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$tr10 = $s10.Shapes.Item(3).TextFrame.TextRange
Set-ParaText $tr10 8 "Connectivity loss (aka liveness heart-beat failure detection)"
